# Completed UI for the job status report
# Adds a new "Job Status" worksheet (after "Tracks") summarising job runs,
# and removes the tab-selection from the previously-active "Artists" sheet
# (Excel moves the selected tab to the newly active sheet automatically).

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last tab ------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Job Status"

# --- Headers --------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Parameters"
$ws.Range("C1").Value = "Start"
$ws.Range("D1").Value = "End"
$ws.Range("E1").Value = "Error"

$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("C1:D1").NumberFormat = "dd\-mm\-yyyy\ hh:mm:ss"

# --- Data rows --------------------------------------------------------------
$ws.Range("A2").Value = "Catalogue Export"
$ws.Range("B2").Value = "JobName = Catalogue Export, FileName = Example.csv"
$ws.Range("C2").Value = 45228.5255921875
$ws.Range("D2").Value = 45228.525594918981

$ws.Range("A3").Value = "Catalogue Export"
$ws.Range("B3").Value = "JobName = Catalogue Export, FileName = Example.csv"
$ws.Range("C3").Value = 45228.530457615743
$ws.Range("D3").Value = 45228.530460266207

$ws.Range("A4").Value = "Catalogue Export"
$ws.Range("B4").Value = "JobName = Catalogue Export, FileName = 2023-10-29 Export.csv"
$ws.Range("C4").Value = 45228.684753900459
$ws.Range("D4").Value = 45228.684756863426

$ws.Range("C2:D4").NumberFormat = "dd\-mm\-yyyy\ hh:mm:ss"

# --- Column widths (approximate AutoFit, matching the original author's
#     manual "best fit" column sizing as closely as this host allows) ------
$ws.Columns("A").ColumnWidth = 13.5
$ws.Columns("B").ColumnWidth = 51
$ws.Columns("C").ColumnWidth = 16.833333333333332
$ws.Columns("D").ColumnWidth = 16.833333333333332
$ws.Columns("E").ColumnWidth = 3.8333333333333335

# --- View: freeze header row, hide gridlines, select next empty row --------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A5").Select() | Out-Null
